$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Re-style the previously-last data row (14) to match the bordered "group"
# look already used by rows 11-13, by copying that formatting over (this
# reuses the existing style/border indices instead of minting new ones).
$ws.Range("A11:E11").Copy()
$ws.Range("A14:E14").PasteSpecial(-4122)  # xlPasteFormats

# Append the new script line as row 15 (English source first, like the
# rest of this sheet's authoring order, then the file/line id, then the
# Russian translation and its "converted"/encoded counterpart).
$ws.Cells.Item(15, 3).Value = ' We\''re out of your league![K]\nYou KO\''d [CS:N]Darkrai[CR]!'
$ws.Cells.Item(15, 1).Value = 'SCRIPT/G01P03A/us2210.ssb'
$ws.Cells.Item(15, 2).Value = 18
$ws.Cells.Item(15, 4).Value = ' Мы явно вам не ровня![K]\nВы уничтожили [CS:N]Даркрая[CR]!'
$ws.Cells.Item(15, 5).Value = ' Íú ÿâîï âàí îå ñïâîÿ![K]\nÂú ôîéœóïçéìé [CS:N]Äàñëñàÿ[CR]!'

# New row keeps the plain (unbordered) look the table used before, sized
# the same as the other wrapped-text rows.
$ws.Rows.Item(15).RowHeight = 43.2

# Restore the saved selection state (active cell on the new last row).
$ws.Range("D15").Select()
